# Add new "LR" results block (rows 65-71) into the "New testing data" table,
# pushing everything below down by 12 rows (matches target dimension A3:Q95).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 12 blank rows right after row 64 (i.e. before the old row 66 block),
# which shifts old rows 66-83 down to 78-95 and keeps rows 62/63 in place.
$ws.Rows("65:76").Insert()

# Extend the Q62/Q63 formula into a shared formula group covering the new rows.
$ws.Range("Q62:Q71").Formula = "=P62*1000*250*200"

# New header row (65) - mirrors row 61's header but for the "LR" sweep.
$ws.Range("B65").Value = "LR"
$ws.Range("C65").Value = "N_CYCLES"
$ws.Range("D65").Value = "WEIGHT"
$ws.Range("E65").Value = "EMB_SIZE"
$ws.Range("F65").Value = "HID_SIZE"
$ws.Range("G65").Value = "ROC_AUC"
$ws.Range("H65").Value = "Avg_prec"
$ws.Range("I65").Value = "Th"
$ws.Range("J65").Value = "TP"
$ws.Range("L65").Value = "FN"
$ws.Range("N65").Value = "FP"
$ws.Range("P65").Value = "TN"

# Row 66
$ws.Range("A66").Value = 1
$ws.Range("B66").Value = 0.0003
$ws.Range("C66").Value = 1
$ws.Range("D66").Value = 1000
$ws.Range("E66").Value = 16
$ws.Range("F66").Value = 32
$ws.Range("G66").Value = 0.958
$ws.Range("H66").Value = 0.0507
$ws.Range("I66").Value = 0.6
$ws.Range("J66").Value = 0.00141
$ws.Range("L66").Value = 0.0002
$ws.Range("N66").Value = 0.111
$ws.Range("P66").Value = 0.887

# Row 67
$ws.Range("A67").Value = 2
$ws.Range("B67").Value = 0.0003
$ws.Range("C67").Value = 1
$ws.Range("D67").Value = 10000
$ws.Range("E67").Value = 16
$ws.Range("F67").Value = 32
$ws.Range("G67").Value = 0.954
$ws.Range("H67").Value = 0.05
$ws.Range("I67").Value = 0.6
$ws.Range("J67").Value = 0.00149
$ws.Range("L67").Value = 0.000122
$ws.Range("N67").Value = 0.151
$ws.Range("P67").Value = 0.847

# Row 68
$ws.Range("A68").Value = 3
$ws.Range("B68").Value = 0.0003
$ws.Range("C68").Value = 1
$ws.Range("D68").Value = 100000
$ws.Range("E68").Value = 16
$ws.Range("F68").Value = 32
$ws.Range("G68").Value = 0.962
$ws.Range("H68").Value = 0.0588
$ws.Range("I68").Value = 0.6
$ws.Range("J68").Value = 0.00114
$ws.Range("L68").Value = 0.000555
$ws.Range("N68").Value = 0.0212
$ws.Range("P68").Value = 0.977

# Row 69
$ws.Range("A69").Value = 4
$ws.Range("B69").Value = 0.0003
$ws.Range("C69").Value = 1
$ws.Range("D69").Value = 1000000
$ws.Range("E69").Value = 16
$ws.Range("F69").Value = 32
$ws.Range("G69").Value = 0.958
$ws.Range("H69").Value = 0.0594
$ws.Range("I69").Value = 0.6
$ws.Range("J69").Value = 0.000648
$ws.Range("L69").Value = 0.000959
$ws.Range("N69").Value = 0.00858
$ws.Range("P69").Value = 0.99

# Row 70
$ws.Range("A70").Value = 5
$ws.Range("B70").Value = 0.0003
$ws.Range("C70").Value = 1
$ws.Range("D70").Value = 50000000
$ws.Range("E70").Value = 16
$ws.Range("F70").Value = 32
$ws.Range("G70").Value = 0.952
$ws.Range("H70").Value = 0.0537
$ws.Range("I70").Value = 0.6
$ws.Range("J70").Value = 0.00117
$ws.Range("L70").Value = 0.000439
$ws.Range("N70").Value = 0.0356
$ws.Range("P70").Value = 0.963

# Row 71
$ws.Range("A71").Value = 6
$ws.Range("B71").Value = 0.0003
$ws.Range("C71").Value = 4
$ws.Range("D71").Value = 2000000
$ws.Range("E71").Value = 16
$ws.Range("F71").Value = 32
$ws.Range("G71").Value = 0.962
$ws.Range("H71").Value = 0.0571
$ws.Range("I71").Value = 0.6
$ws.Range("J71").Value = 0.00155
$ws.Range("L71").Value = 0.0000523
$ws.Range("N71").Value = 0.18
$ws.Range("P71").Value = 0.818

# Shared formulas for the new K/M/O columns (rows 66-71).
$ws.Range("K66:K71").Formula = "=J66*1000*250*200"
$ws.Range("M66:M71").Formula = "=L66*1000*250*200"
$ws.Range("O66:O71").Formula = "=N66*1000*250*200"

Write-Host "edit complete"
